$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "614×8=" "811×8="
Replace-Text "359×8=" "483×9="
Replace-Text "848×3=" "383×9="
Replace-Text "750×7=" "895×9="
Replace-Text "766×9=" "146×6="
Replace-Text "344×6=" "251×8="
Replace-Text "991×2=" "827×4="
Replace-Text "837×2=" "923×5="
Replace-Text "294×5=" "498×5="
Replace-Text "904×5=" "914×8="
Replace-Text "255×9=" "823×2="
Replace-Text "151×6=" "873×5="
Replace-Text "921×8=" "466×8="
Replace-Text "953×2=" "383×7="
Replace-Text "778×2=" "841×8="
Replace-Text "137×9=" "514×4="
Replace-Text "473×3=" "628×6="
Replace-Text "166×3=" "816×2="
Replace-Text "849×4=" "629×4="
Replace-Text "506×2=" "499×9="
Replace-Text "266×5=" "221×9="
Replace-Text "124×8=" "249×3="
Replace-Text "166×4=" "983×7="
Replace-Text "907×9=" "686×5="
Replace-Text "309×8=" "541×3="
